$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the mileage log with the entries recorded during lockdown ---
# Column A: one row per day from 2020-03-09 (serial 43899) to 2020-04-09 (serial 43930)
# Column B: odometer reading for that day

$firstRow = 45
$lastRow = 76
$firstDateSerial = 43899

for ($i = 0; $i -le ($lastRow - $firstRow); $i++) {
    $r = $firstRow + $i
    $ws.Cells.Item($r, 1).Value = $firstDateSerial + $i
}

# B45:B74 all use the same relative formula (value = next day's reading - 1),
# set on the whole block at once so Excel stores it as one shared formula.
$ws.Range("B45:B74").Formula = "=B46-1"

# B75 keeps the same kind of formula but Excel no longer folds it into the
# B45:B74 shared-formula group because it is the second-to-last row.
$ws.Range("B75").Formula = "=B76-1"

# B76 is the last, literal odometer reading for the period.
$ws.Range("B76").Value = 4713

# --- Restore the view the author left the sheet in ---
$ws.Range("C45").Select()
$excel.ActiveWindow.ScrollRow = 41
$excel.ActiveWindow.ScrollColumn = 1
